$d = $word.ActiveDocument

# --- Paragraph 1 (index 1): title line - two runs separated by a manual line break ---
$p1 = $d.Paragraphs.Item(1)
$r1 = $p1.Range
$found1 = $r1.Find.Execute("המאמר היומי של מייק - 03.03.25", $true, $false, $false, $false, $false, $true, 1, $false, @'
המאמר היומי של מייק - 02.03.25
'@, 2)
if (-not $found1) { throw "Could not find the review-date title text to replace" }

$r1b = $p1.Range
$found2 = $r1b.Find.Execute(" The Geometry of Concepts: Sparse Autoencoder Feature Structure", $true, $false, $false, $false, $false, $true, 1, $false, @'
An Empirical Model of Large-Batch Training
'@, 2)
if (-not $found2) { throw "Could not find the paper-title text to replace" }

# --- Helper: replace the body text of a whole paragraph (keeps pPr, single run) ---
function Set-ParagraphText($para, $newText) {
    $full = $para.Range
    $body = $d.Range($full.Start, $full.End - 1)
    $body.Text = $newText
}

Set-ParagraphText $d.Paragraphs.Item(2) @'
מאמר מלפני 6 שנים של חוקרי OpenAI אך מצאתי אותו די מעניין לסקירה קצרה. המאמר חוקר גדול באץ' אופטימלי עבור אימון Mini-Batch Gradient Descent או MBGD. מה זה אופטימלי כאן? כזה שימעזר את מספר הדוגמאות ש-MiGD משתמש בהם כדי להביא את המודל לערך יעד של הלוס. כמובן שניתן ״להריץ״ את אותו הדוגמא כמה פעמים במהלך MBGD.
'@

Set-ParagraphText $d.Paragraphs.Item(3) @'
למי שכח MBGD שייך למשפחת שיטות המבוססות על מורד הגרדיאנט. עם MBGD אנו מחלקים את הדאטהסט למיני-באצ'ים שכל באץ' מורכב מכמה דוגמאות. עבור כל באץ' אנו מבצעים עדכון אחד של משקלי מודל כאשר הגרדיאנט מחושב בתור ממוצע של כל ערכי הגרדיאנטים עבור כל הדוגמאות בבאץ'. למעשה ממוצע זה הינו משערך של הגרדיאנט הממוצע של המודל עבור כל הדוגמאות מהדאטהסט. נזכיר שכל עדכון הוא הזזה (לינארית) של משקולות המודל בכיוון ההפוך לכיוון הגרדיאנט. כל עדכון כזה תלוי בקצב למידה שקובע את גודל עדכון המשקולות (מוכפל בגרדיאנט ממוצע).
'@

Set-ParagraphText $d.Paragraphs.Item(4) @'
המאמר מציע שיטה למציאת גודל באץ' אופטימלי (לפי ההגדרה שנתתי קודם) שעבור קצב למידה אופטימלי (הממזער את הלוס בכל איטרציה). די ברור כי גודל באץ' אופטימלי צריך להיות תלוי בפרמטרי המודל - למשל בצורת משטח הלוס וגם בערכי הגרדיאנט. המאמר טוען כי גודל באץ' אופטימלי ניתן לחשב בתור הטרייס (trace, סכום הערכים העצמיים) של המכפלה של מטריצת קווריאנס של גרדיאנט הלוס וההיסאין H של פונקציית לוס מחולקת ב G^T)HG) כאשר G הוא הממוצע של וקטור הגרדיאנט.
'@

Set-ParagraphText $d.Paragraphs.Item(5) @'
תוצאה זו התקבלה דרך פיתוח טיילור מסדר שני (בכיוון הגרדיאנט)ֿ, מציאה גודל קצב למידה אופטימלי והצבתו לנוסחה כדי לחשב את גודל הבאץ' שעבורו מתקבל ירידה מקסימלית של הלוס. לאחר מכן משווים את הירידה המקסימלית  עם זו עבור גודל באץ' נתון B.
'@

Set-ParagraphText $d.Paragraphs.Item(6) @'
המאמר מדגיש שגודל אופטימלי של באץ' אינו תלוי בגודל שלה דאטהסט וכמובן משתנה במהלך האימון כי גם ההסיאן H וגם הגרדיאנט הממוצע H וגם מטריצת קווריאנס של גרדיאנט הלוס לא נשארים קבועים (בד״כ). המחברים מציינים מקרה פרטי די מעניין (לא קורה במציאות אמנם) שבו ההיסאין H שווה למטריצה היחידה I. במקרה הזה גודל באץ' אופטימלי שווה לסכום השונויות של כל רכיבי הגרדיאנט.
'@

Set-ParagraphText $d.Paragraphs.Item(7) @'
המאמר כתוב בצורה מאוד מובנת וניתן לקריאה קלילה יחסית…
'@

# --- Remove the five paragraphs (8-12) describing the old "atomic/brain/galactic" sections ---
$countBefore = $d.Paragraphs.Count
$pStart = $d.Paragraphs.Item(8)
$pEnd = $d.Paragraphs.Item(12)
$delRange = $d.Range($pStart.Range.Start, $pEnd.Range.End)
$delRange.Delete()
$countAfter = $d.Paragraphs.Count
if ($countAfter -ne ($countBefore - 5)) { throw "Unexpected paragraph count after delete: before=$countBefore after=$countAfter" }

# --- Last paragraph: update the arxiv link ---
Set-ParagraphText $d.Paragraphs.Item(8) @'
https://arxiv.org/abs/1812.06162
'@

Write-Output "Final paragraph count: $($d.Paragraphs.Count)"
